$wb = $excel.ActiveWorkbook

# Update 2022 (column I) violent crime totals with data through 2022-07-11
# Each block below targets one worksheet and sets the new column-I values
# for the rows whose underlying crime counts changed.

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I6").Value = 22  # Robbery: 21 -> 22
$ws.Range("I7").Value = 112  # Total: 111 -> 112

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I6").Value = 107  # Robbery: 106 -> 107
$ws.Range("I7").Value = 414  # Total: 413 -> 414

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I3").Value = 215  # Aggravated Battery: 214 -> 215
$ws.Range("I6").Value = 250  # Robbery: 247 -> 250
$ws.Range("I7").Value = 775  # Total: 771 -> 775

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I2").Value = 21  # Aggravated Assault: 20 -> 21
$ws.Range("I7").Value = 58  # Total: 57 -> 58

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 28  # Aggravated Assault: 27 -> 28
$ws.Range("I7").Value = 88  # Total: 87 -> 88

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I3").Value = 17  # Aggravated Battery: 16 -> 17
$ws.Range("I7").Value = 64  # Total: 63 -> 64

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("I6").Value = 23  # Robbery: 20 -> 23
$ws.Range("I7").Value = 37  # Total: 34 -> 37

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 112  # Albany Park: 111 -> 112
$ws.Range("I7").Value = 414  # Auburn Gresham: 413 -> 414
$ws.Range("I8").Value = 775  # Austin: 771 -> 775
$ws.Range("I9").Value = 58  # Avalon Park: 57 -> 58
$ws.Range("I10").Value = 88  # Avondale: 87 -> 88
$ws.Range("I14").Value = 64  # Bridgeport: 63 -> 64
$ws.Range("I16").Value = 37  # Bucktown: 34 -> 37
$ws.Range("I18").Value = 87  # Calumet Heights: 85 -> 87
$ws.Range("I19").Value = 340  # Chatham: 336 -> 340
$ws.Range("I20").Value = 315  # Chicago Lawn: 314 -> 315
$ws.Range("I23").Value = 121  # Douglas: 120 -> 121
$ws.Range("I25").Value = 61  # East Side: 60 -> 61
$ws.Range("I27").Value = 116  # Edgewater: 114 -> 116
$ws.Range("I29").Value = 831  # Englewood: 829 -> 831
$ws.Range("I33").Value = 576  # Garfield Park: 575 -> 576
$ws.Range("I36").Value = 176  # Grand Boulevard: 175 -> 176
$ws.Range("I37").Value = 408  # Grand Crossing: 406 -> 408
$ws.Range("I41").Value = 57  # Hermosa: 56 -> 57
$ws.Range("I42").Value = 446  # Humboldt Park: 444 -> 446
$ws.Range("I48").Value = 165  # Lake View: 164 -> 165
$ws.Range("I51").Value = 123  # Little Italy, UIC: 122 -> 123
$ws.Range("I52").Value = 278  # Little Village: 277 -> 278
$ws.Range("I54").Value = 285  # Loop: 283 -> 285
$ws.Range("I59").Value = 26  # Montclare: 25 -> 26
$ws.Range("I63").Value = 50  # NO NEIGHBORHOOD DATA: 46 -> 50
$ws.Range("I65").Value = 282  # New City: 280 -> 282
$ws.Range("I67").Value = 492  # North Lawndale: 486 -> 492
$ws.Range("I68").Value = 41  # North Park: 40 -> 41
$ws.Range("I69").Value = 29  # Norwood Park: 28 -> 29
$ws.Range("I78").Value = 184  # Rogers Park: 182 -> 184
$ws.Range("I79").Value = 338  # Roseland: 333 -> 338
$ws.Range("I83").Value = 257  # South Chicago: 252 -> 257
$ws.Range("I85").Value = 581  # South Shore: 580 -> 581
$ws.Range("I89").Value = 143  # Uptown: 142 -> 143
$ws.Range("I91").Value = 154  # Washington Park: 152 -> 154
$ws.Range("I92").Value = 38  # West Elsdon: 37 -> 38
$ws.Range("I97").Value = 91  # West Town: 90 -> 91
$ws.Range("I99").Value = 238  # Woodlawn: 237 -> 238
$ws.Range("I101").Value = 12774  # Total: 12703 -> 12774

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I2").Value = 24  # Aggravated Assault: 23 -> 24
$ws.Range("I3").Value = 23  # Aggravated Battery: 22 -> 23
$ws.Range("I7").Value = 87  # Total: 85 -> 87

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 130  # Aggravated Assault: 129 -> 130
$ws.Range("I3").Value = 97  # Aggravated Battery: 95 -> 97
$ws.Range("I6").Value = 93  # Robbery: 92 -> 93
$ws.Range("I7").Value = 340  # Total: 336 -> 340

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I3").Value = 96  # Aggravated Battery: 95 -> 96
$ws.Range("I7").Value = 315  # Total: 314 -> 315

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 3610  # Aggravated Assault: 3591 -> 3610
$ws.Range("I3").Value = 3742  # Aggravated Battery: 3725 -> 3742
$ws.Range("I4").Value = 873  # Criminal Sexual Assault: 872 -> 873
$ws.Range("I5").Value = 347  # Homicide: 343 -> 347
$ws.Range("I6").Value = 4202  # Robbery: 4172 -> 4202
$ws.Range("I7").Value = 12774  # Total: 12703 -> 12774

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I6").Value = 37  # Robbery: 36 -> 37
$ws.Range("I7").Value = 121  # Total: 120 -> 121

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I3").Value = 17  # Aggravated Battery: 16 -> 17
$ws.Range("I7").Value = 61  # Total: 60 -> 61

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 30  # Aggravated Assault: 29 -> 30
$ws.Range("I6").Value = 47  # Robbery: 46 -> 47
$ws.Range("I7").Value = 116  # Total: 114 -> 116

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I5").Value = 31  # Homicide: 30 -> 31
$ws.Range("I6").Value = 227  # Robbery: 226 -> 227
$ws.Range("I7").Value = 831  # Total: 829 -> 831

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I6").Value = 184  # Robbery: 183 -> 184
$ws.Range("I7").Value = 576  # Total: 575 -> 576

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I5").Value = 7  # Homicide: 6 -> 7
$ws.Range("I7").Value = 176  # Total: 175 -> 176

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 131  # Aggravated Assault: 130 -> 131
$ws.Range("I3").Value = 124  # Aggravated Battery: 123 -> 124
$ws.Range("I7").Value = 408  # Total: 406 -> 408

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I2").Value = 20  # Aggravated Assault: 19 -> 20
$ws.Range("I7").Value = 57  # Total: 56 -> 57

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 152  # Aggravated Battery: 151 -> 152
$ws.Range("I6").Value = 120  # Robbery: 119 -> 120
$ws.Range("I7").Value = 446  # Total: 444 -> 446

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I4").Value = 16  # Criminal Sexual Assault: 17 -> 16
$ws.Range("I6").Value = 93  # Robbery: 91 -> 93
$ws.Range("I7").Value = 165  # Total: 164 -> 165

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 49  # Robbery: 48 -> 49
$ws.Range("I7").Value = 123  # Total: 122 -> 123

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I6").Value = 67  # Robbery: 66 -> 67
$ws.Range("I7").Value = 278  # Total: 277 -> 278

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I3").Value = 36  # Aggravated Battery: 37 -> 36
$ws.Range("I6").Value = 64  # Robbery: 63 -> 64

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 56  # Aggravated Battery: 55 -> 56
$ws.Range("I6").Value = 145  # Robbery: 144 -> 145
$ws.Range("I7").Value = 285  # Total: 283 -> 285

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("I3").Value = 4  # Aggravated Battery: 3 -> 4
$ws.Range("I7").Value = 26  # Total: 25 -> 26

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 29  # Aggravated Assault: 28 -> 29
$ws.Range("I6").Value = 42  # Robbery: 43 -> 42

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 77  # Aggravated Battery: 76 -> 77
$ws.Range("I5").Value = 14  # Homicide: 13 -> 14
$ws.Range("I7").Value = 282  # Total: 280 -> 282

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 118  # Aggravated Assault: 117 -> 118
$ws.Range("I3").Value = 172  # Aggravated Battery: 170 -> 172
$ws.Range("I4").Value = 25  # Criminal Sexual Assault: 24 -> 25
$ws.Range("I5").Value = 12  # Homicide: 11 -> 12
$ws.Range("I6").Value = 165  # Robbery: 164 -> 165
$ws.Range("I7").Value = 492  # Total: 486 -> 492

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("I6").Value = 8  # Robbery: 7 -> 8
$ws.Range("I7").Value = 41  # Total: 40 -> 41

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("I6").Value = 10  # Robbery: 9 -> 10
$ws.Range("I7").Value = 29  # Total: 28 -> 29

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I2").Value = 37  # Aggravated Assault: 36 -> 37
$ws.Range("I3").Value = 47  # Aggravated Battery: 46 -> 47
$ws.Range("I7").Value = 184  # Total: 182 -> 184

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 99  # Aggravated Assault: 97 -> 99
$ws.Range("I3").Value = 110  # Aggravated Battery: 109 -> 110
$ws.Range("I6").Value = 101  # Robbery: 99 -> 101
$ws.Range("I7").Value = 338  # Total: 333 -> 338

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 89  # Aggravated Assault: 87 -> 89
$ws.Range("I3").Value = 100  # Aggravated Battery: 98 -> 100
$ws.Range("I6").Value = 47  # Robbery: 46 -> 47
$ws.Range("I7").Value = 257  # Total: 252 -> 257

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 151  # Aggravated Assault: 150 -> 151
$ws.Range("I3").Value = 232  # Aggravated Battery: 233 -> 232
$ws.Range("I6").Value = 146  # Robbery: 145 -> 146
$ws.Range("I7").Value = 581  # Total: 580 -> 581

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 35  # Aggravated Assault: 34 -> 35
$ws.Range("I7").Value = 143  # Total: 142 -> 143

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I2").Value = 49  # Aggravated Assault: 48 -> 49
$ws.Range("I6").Value = 48  # Robbery: 47 -> 48
$ws.Range("I7").Value = 154  # Total: 152 -> 154

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I2").Value = 14  # Aggravated Assault: 13 -> 14
$ws.Range("I7").Value = 38  # Total: 37 -> 38

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I3").Value = 17  # Aggravated Battery: 16 -> 17
$ws.Range("I7").Value = 91  # Total: 90 -> 91

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I6").Value = 65  # Robbery: 64 -> 65
$ws.Range("I7").Value = 238  # Total: 237 -> 238
